# Append a new blank/zeroed data row (row 2) below the header row on the
# active sheet ("Đơn sale chính"), matching the report-format update.
#
# Text / label columns -> blank ("")
# Money / numeric columns -> 0

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text columns: A, C, D, E, F, G, H, J, Q, R, S, T
$textCols = @("A2","C2","D2","E2","F2","G2","H2","J2","Q2","R2","S2","T2")
foreach ($addr in $textCols) {
    $ws.Range($addr).Formula = '=""'
}

# Numeric columns: B, I, K, L, M, N, O, P
$numCols = @("B2","I2","K2","L2","M2","N2","O2","P2")
foreach ($addr in $numCols) {
    $ws.Range($addr).Value = 0
}
